# "tidy set definitions up a little ... changed input headers to all start from 0"
#
# Column J held an input "header" value of 100000000 in J6, which every row
# below (J7:J29) picks up via a shared formula "=<cell above>". Updating J6
# propagates the new value down the whole shared-formula chain.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the source input value; J7:J29 (=J6, =J7, ... shared formulas)
# recalculate automatically from this.
$ws.Range("J6").Value = 9999

# Recalculate the workbook so all dependent cached formula values are refreshed.
$excel.CalculateFullRebuild()

# Restore/update the sheet's selection state: top-right pane back to its
# normal default cell, and the active/frozen bottom-right pane pointed at J7
# (single cell) instead of the old Z6:Z29 selection.
$ws.Range("B1").Select()
$ws.Range("J7").Select()
